# Apply the "Updated column width calclation" edit to the "Columns" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Columns")

# Swap the Column Width values for the "Title" (row 4) and "Sex" (row 5) rows.
$ws.Range("D4").Value = 24
$ws.Range("D5").Value = 12

# Fill in the missing Text Align / Font Bold / Font Name values for the
# "Eval Number" row (row 8), matching the pattern used by the other rows.
$ws.Range("G8").Value = "Left"
$ws.Range("H8").Value = "N"
$ws.Range("I8").Value = "Calibri"

# Update the selection to match the cells that were just edited.
$ws.Range("H8:I8").Select()
